# Update "想去人数" (interest count) values in column F across sheets,
# matching the author's regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(4, 6).Value = 105
$ws.Cells.Item(5, 6).Value = 69
$ws.Cells.Item(6, 6).Value = 708
$ws.Cells.Item(8, 6).Value = 206
$ws.Cells.Item(9, 6).Value = 16
$ws.Cells.Item(11, 6).Value = 35
$ws.Cells.Item(12, 6).Value = 568
$ws.Cells.Item(16, 6).Value = 139
$ws.Cells.Item(17, 6).Value = 779
$ws.Cells.Item(18, 6).Value = 2581
$ws.Cells.Item(21, 6).Value = 7
$ws.Cells.Item(22, 6).Value = 303
$ws.Cells.Item(25, 6).Value = 135
$ws.Cells.Item(27, 6).Value = 947
$ws.Cells.Item(29, 6).Value = 163

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 537
$ws.Cells.Item(5, 6).Value = 537
$ws.Cells.Item(7, 6).Value = 14
$ws.Cells.Item(8, 6).Value = 220
$ws.Cells.Item(14, 6).Value = 533
$ws.Cells.Item(15, 6).Value = 82
$ws.Cells.Item(17, 6).Value = 956
$ws.Cells.Item(22, 6).Value = 32
$ws.Cells.Item(24, 6).Value = 268
$ws.Cells.Item(31, 6).Value = 24

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 2356
$ws.Cells.Item(10, 6).Value = 306
$ws.Cells.Item(11, 6).Value = 86

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value = 2356
$ws.Cells.Item(10, 6).Value = 306
$ws.Cells.Item(11, 6).Value = 86
$ws.Cells.Item(12, 6).Value = 105
$ws.Cells.Item(13, 6).Value = 69
$ws.Cells.Item(14, 6).Value = 708
$ws.Cells.Item(17, 6).Value = 206
$ws.Cells.Item(18, 6).Value = 16
$ws.Cells.Item(19, 6).Value = 35
$ws.Cells.Item(20, 6).Value = 568
$ws.Cells.Item(21, 6).Value = 537
$ws.Cells.Item(24, 6).Value = 14
$ws.Cells.Item(25, 6).Value = 139
$ws.Cells.Item(26, 6).Value = 779
$ws.Cells.Item(27, 6).Value = 2581
$ws.Cells.Item(30, 6).Value = 303
$ws.Cells.Item(32, 6).Value = 135
$ws.Cells.Item(34, 6).Value = 947
$ws.Cells.Item(35, 6).Value = 533
$ws.Cells.Item(36, 6).Value = 82
$ws.Cells.Item(38, 6).Value = 163
$ws.Cells.Item(42, 6).Value = 32
$ws.Cells.Item(43, 6).Value = 268
$ws.Cells.Item(44, 6).Value = 268
$ws.Cells.Item(48, 6).Value = 24
